# Daily attendance processing - 2025-12-09 19:48:15
# Normalizes the "Recorded By" (column G) cell values: for any cell whose
# value is a comma-separated list of recorders, swap the first and last
# entries in the list (leaving any middle entries untouched).
#
# Rows 4, 30 and 56 (the "session 3" row for each of the B2A/B2B/B2C
# blocks, all holding "System, backup@backdoor.com") were already
# normalized in an earlier pass today and are left untouched here.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

$skipRows = @(4, 30, 56)

for ($r = 2; $r -le $lastRow; $r++) {
    if ($skipRows -contains $r) { continue }

    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $val = $cell.Value2

    if ($val -eq $null) { continue }

    $parts = $val -split ', '
    if ($parts.Count -gt 1) {
        $first = $parts[0]
        $last = $parts[$parts.Count - 1]
        $parts[0] = $last
        $parts[$parts.Count - 1] = $first
        $cell.Value2 = [string]::Join(', ', $parts)
    }
}
